# I0 and IF added
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting (bold, bordered, centered) from the existing
# "IP" header (H1) onto the two new header cells so they pick up the same
# cell style used by the other column headers.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data columns
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1

$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 7

$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 6

$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 5

$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 6

$ws.Range("I7").Value = 1
$ws.Range("J7").Value = 6

$ws.Range("I8").Value = 1
$ws.Range("J8").Value = 6

$ws.Range("I9").Value = 1
$ws.Range("J9").Value = 6

$ws.Range("I10").Value = 1
$ws.Range("J10").Value = 7

$ws.Range("I11").Value = 1
$ws.Range("J11").Value = 6

$ws.Range("I12").Value = 1
$ws.Range("J12").Value = 4

$ws.Range("I13").Value = 1
$ws.Range("J13").Value = 4

$ws.Range("I14").Value = 4
$ws.Range("J14").Value = 6

$ws.Range("I15").Value = 1
$ws.Range("J15").Value = 2
